$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1991.9584
$ws.Range("I116").Value = 1567.5
$ws.Range("J116").Value = 2133.4443
$ws.Range("K116").Value = 1567.5
$ws.Range("L116").Value = 2133.4443
$ws.Range("M116").Value = 1874.5
$ws.Range("N116").Value = -9017.444299999999

$ws.Range("H137").Value = 1498.8667
$ws.Range("I137").Value = 1267.9231
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 3803.7693
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -1253.7693
$ws.Range("N137").Value = -14100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17054.816
$ws.Range("I32").Value = 18257.469
$ws.Range("K32").Value = 18257.469
$ws.Range("M32").Value = -17970.469

$ws.Range("H74").Value = 664.6389
$ws.Range("I74").Value = 538.69696
$ws.Range("K74").Value = 538.69696
$ws.Range("M74").Value = 335.30304

$ws.Range("H77").Value = 664.6389
$ws.Range("I77").Value = 538.69696
$ws.Range("K77").Value = 2693.4848
$ws.Range("M77").Value = 1674.5152

$ws.Range("H97").Value = 621.8929000000001
$ws.Range("I97").Value = 642.65
$ws.Range("J97").Value = 570
$ws.Range("K97").Value = 642.65
$ws.Range("L97").Value = 570
$ws.Range("M97").Value = -146.65
$ws.Range("N97").Value = -1562

$ws.Range("H122").Value = 1413.9474
$ws.Range("I122").Value = 1409.1666
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4227.4998
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1777.4998
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 3639.282
$ws.Range("I132").Value = 3981.0667
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 11943.2001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -9413.2001
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 540
$ws.Range("I80").Value = 925
$ws.Range("J80").Value = 386
$ws.Range("K80").Value = 925
$ws.Range("L80").Value = 386
$ws.Range("M80").Value = 73
$ws.Range("N80").Value = -2382

$ws.Range("H83").Value = 540
$ws.Range("I83").Value = 925
$ws.Range("J83").Value = 386
$ws.Range("K83").Value = 4625
$ws.Range("L83").Value = 1930
$ws.Range("M83").Value = 367
$ws.Range("N83").Value = -11914

$ws.Range("H105").Value = 1877.2972
$ws.Range("I105").Value = 1717.3914
$ws.Range("J105").Value = 2140
$ws.Range("K105").Value = 1717.3914
$ws.Range("L105").Value = 2140
$ws.Range("M105").Value = 29.60860000000002
$ws.Range("N105").Value = -5634

$ws.Range("H134").Value = 36680.9
$ws.Range("I134").Value = 49483.19
$ws.Range("J134").Value = 3074.875
$ws.Range("K134").Value = 148449.57
$ws.Range("L134").Value = 9224.625
$ws.Range("M134").Value = -145914.57
$ws.Range("N134").Value = -14294.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.28570999999999
$ws.Range("I7").Value = 70.72727
$ws.Range("J7").Value = 176
$ws.Range("K7").Value = 70.72727
$ws.Range("L7").Value = 176
$ws.Range("M7").Value = 42.27273
$ws.Range("N7").Value = -402

$ws.Range("H16").Value = 2493.3333
$ws.Range("I16").Value = 980
$ws.Range("J16").Value = 2796
$ws.Range("K16").Value = 980
$ws.Range("L16").Value = 2796
$ws.Range("M16").Value = -693
$ws.Range("N16").Value = -3370

$ws.Range("H31").Value = 9095926
$ws.Range("I31").Value = 3354.611
$ws.Range("K31").Value = 3354.611
$ws.Range("M31").Value = -3059.611

$ws.Range("H34").Value = 9095926
$ws.Range("I34").Value = 3354.611
$ws.Range("K34").Value = 3354.611
$ws.Range("M34").Value = -3152.611

$ws.Range("H58").Value = 842.2162
$ws.Range("I58").Value = 798.4706
$ws.Range("K58").Value = 798.4706
$ws.Range("M58").Value = -595.4706

$ws.Range("H94").Value = 76926630
$ws.Range("I94").Value = 142859970
$ws.Range("J94").Value = 4416.6665
$ws.Range("K94").Value = 142859970
$ws.Range("L94").Value = 4416.6665
$ws.Range("M94").Value = -142859519
$ws.Range("N94").Value = -5318.6665

$ws.Range("H99").Value = 1370.9678
$ws.Range("I99").Value = 900
$ws.Range("K99").Value = 900
$ws.Range("M99").Value = 598

$ws.Range("H105").Value = 1132
$ws.Range("I105").Value = 551.6667
$ws.Range("J105").Value = 2002.5
$ws.Range("K105").Value = 551.6667
$ws.Range("L105").Value = 2002.5
$ws.Range("M105").Value = 1195.3333
$ws.Range("N105").Value = -5496.5

$ws.Range("H107").Value = 470.91306
$ws.Range("I107").Value = 512.375
$ws.Range("J107").Value = 376.14285
$ws.Range("K107").Value = 512.375
$ws.Range("L107").Value = 376.14285
$ws.Range("M107").Value = 1407.625
$ws.Range("N107").Value = -4216.14285

$ws.Range("H113").Value = 2493.3333
$ws.Range("I113").Value = 980
$ws.Range("J113").Value = 2796
$ws.Range("K113").Value = 980
$ws.Range("L113").Value = 2796
$ws.Range("M113").Value = 1190
$ws.Range("N113").Value = -7136

$ws.Range("H122").Value = 620.34485
$ws.Range("I122").Value = 565.5238000000001
$ws.Range("J122").Value = 764.25
$ws.Range("K122").Value = 1696.5714
$ws.Range("L122").Value = 2292.75
$ws.Range("M122").Value = 753.4285999999997
$ws.Range("N122").Value = -7192.75

$ws.Range("H126").Value = 1370.9678
$ws.Range("I126").Value = 900
$ws.Range("K126").Value = 2700
$ws.Range("M126").Value = -230

$ws.Range("H132").Value = 2500.182
$ws.Range("I132").Value = 1706.2354
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 5118.706200000001
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -2588.706200000001
$ws.Range("N132").Value = -20658.8

$ws.Range("H134").Value = 1113.1
$ws.Range("I134").Value = 1014.55554
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3043.66662
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -508.66662
$ws.Range("N134").Value = -11070

$ws.Range("H136").Value = 842.2162
$ws.Range("I136").Value = 798.4706
$ws.Range("K136").Value = 2395.4118
$ws.Range("M136").Value = 154.5882000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2501
$ws.Range("J48").Value = 2501
$ws.Range("L48").Value = 7503
$ws.Range("N48").Value = -8003

$ws.Range("H63").Value = 8596
$ws.Range("I63").Value = 5490
$ws.Range("K63").Value = 16470
$ws.Range("M63").Value = -15721

$ws.Range("H64").Value = 1858.5714
$ws.Range("I64").Value = 973
$ws.Range("J64").Value = 2100.0908
$ws.Range("K64").Value = 2919
$ws.Range("L64").Value = 6300.2724
$ws.Range("M64").Value = -2649
$ws.Range("N64").Value = -6840.2724

$ws.Range("H66").Value = 8596
$ws.Range("I66").Value = 5490
$ws.Range("K66").Value = 49410
$ws.Range("M66").Value = -45666

$ws.Range("H67").Value = 1858.5714
$ws.Range("I67").Value = 973
$ws.Range("J67").Value = 2100.0908
$ws.Range("K67").Value = 2919
$ws.Range("L67").Value = 6300.2724
$ws.Range("M67").Value = -1983
$ws.Range("N67").Value = -8172.2724

$ws.Range("H70").Value = 4071.0715
$ws.Range("I70").Value = 1997.5
$ws.Range("J70").Value = 4416.6665
$ws.Range("K70").Value = 5992.5
$ws.Range("L70").Value = 13249.9995
$ws.Range("M70").Value = -5677.5
$ws.Range("N70").Value = -13879.9995

$ws.Range("H73").Value = 4071.0715
$ws.Range("I73").Value = 1997.5
$ws.Range("J73").Value = 4416.6665
$ws.Range("K73").Value = 5992.5
$ws.Range("L73").Value = 13249.9995
$ws.Range("M73").Value = -4900.5
$ws.Range("N73").Value = -15433.9995

$ws.Range("H121").Value = 811.7273
$ws.Range("I121").Value = 462
$ws.Range("J121").Value = 1103.1666
$ws.Range("K121").Value = 1386
$ws.Range("L121").Value = 3309.4998
$ws.Range("M121").Value = -76
$ws.Range("N121").Value = -5929.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7290329.5
$ws.Range("I70").Value = 9627002
$ws.Range("J70").Value = 5410.5293
$ws.Range("K70").Value = 9627002
$ws.Range("L70").Value = 5410.5293
$ws.Range("M70").Value = -9626732
$ws.Range("N70").Value = -5950.5293

$ws.Range("H73").Value = 7290329.5
$ws.Range("I73").Value = 9627002
$ws.Range("J73").Value = 5410.5293
$ws.Range("K73").Value = 9627002
$ws.Range("L73").Value = 5410.5293
$ws.Range("M73").Value = -9626066
$ws.Range("N73").Value = -7282.5293

$ws.Range("H80").Value = 6884.643
$ws.Range("I80").Value = 4936.875
$ws.Range("J80").Value = 9481.666999999999
$ws.Range("K80").Value = 4936.875
$ws.Range("L80").Value = 9481.666999999999
$ws.Range("M80").Value = -3938.875
$ws.Range("N80").Value = -11477.667

$ws.Range("H83").Value = 6884.643
$ws.Range("I83").Value = 4936.875
$ws.Range("J83").Value = 9481.666999999999
$ws.Range("K83").Value = 24684.375
$ws.Range("L83").Value = 47408.335
$ws.Range("M83").Value = -19692.375
$ws.Range("N83").Value = -57392.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1780.1212
$ws.Range("I46").Value = 1516.75
$ws.Range("J46").Value = 1930.619
$ws.Range("K46").Value = 1516.75
$ws.Range("L46").Value = 1930.619
$ws.Range("M46").Value = -1328.75
$ws.Range("N46").Value = -2306.619

$ws.Range("H82").Value = 1432.7273
$ws.Range("I82").Value = 1484.4445
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 1484.4445
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -1123.4445
$ws.Range("N82").Value = -1922

$ws.Range("H85").Value = 1432.7273
$ws.Range("I85").Value = 1484.4445
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 1484.4445
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -236.4445000000001
$ws.Range("N85").Value = -3696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 435
$ws.Range("I100").Value = 442
$ws.Range("K100").Value = 884
$ws.Range("M100").Value = -343

$ws.Range("H122").Value = 1149.909
$ws.Range("I122").Value = 1117.6666
$ws.Range("J122").Value = 1295
$ws.Range("K122").Value = 3352.9998
$ws.Range("L122").Value = 3885
$ws.Range("M122").Value = -902.9998000000001
$ws.Range("N122").Value = -8785

$ws.Range("H136").Value = 7554.421
$ws.Range("I136").Value = 8514.625
$ws.Range("J136").Value = 2433.3333
$ws.Range("K136").Value = 25543.875
$ws.Range("L136").Value = 7299.999899999999
$ws.Range("M136").Value = -22993.875
$ws.Range("N136").Value = -12399.9999
